$d = $word.ActiveDocument

# 1. Remove the stray empty paragraph that follows the two-tab paragraph
#    (just before "Unemployment and Full Employment").
$found = $d.Content.Find.Execute("Unemployment and Full Employment", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $target = $d.Content.Find.Parent
    $heading = $target.Paragraphs.Item(1)
    $prev = $heading.Previous()
    if ($prev.Range.Text.Trim().Length -eq 0) {
        $prev.Range.Delete()
    }
}

# 2. Move the "_GoBack" bookmark out of the empty paragraph before "Core Inflation Rate"
#    and into the "core CPI inflation rate" paragraph, splitting the sentence exactly
#    where the new text run boundary goes ("...is c" | "ore CPI inflation rate.").
$gb = $d.Bookmarks.Item("_GoBack")
$gb.Delete()

$found2 = $d.Content.Find.Execute("is core CPI inflation rate", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $hit = $d.Content.Find.Parent
    $splitPoint = $hit.Start + 4
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
